$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "session_001"
$ws.Range("B21").Value = 7.328892469406128
$ws.Range("C21").Value = 2.722530126571655
$ws.Range("D21").Value = 1.859208822250366
$ws.Range("E21").Value = 11.91063141822815
